# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking "Price" strings are prefixed with a leading apostrophe so
# Excel stores them as text (preserving formatting like trailing zeros /
# thousand-dot groupings) instead of silently coercing them to doubles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.063.24'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '3.108.88'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''577.81'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = '''178.22'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.108.36'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.512'
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('D10').Value = '''6.36'
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = '''0.470'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('E13').Value = '  -2.72%  '
$ws.Range('D14').Value = '''36.24'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '3.626.93'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '66.976.30'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').Value = '''7.01'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').Value = '''16.93'
$ws.Range('E19').Value = '  +2.16%  '
$ws.Range('D20').Value = '3.107.55'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = '''485.66'
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('D22').Value = '''7.75'
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').Value = '''0.691'
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = '''83.71'
$ws.Range('D25').Value = '''12.68'
$ws.Range('E25').Value = '  -3.52%  '
$ws.Range('D26').Value = '''2.24'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('D27').Value = '''10.28'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').Value = '''8.08'
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('D30').Value = '''2.29'
$ws.Range('E30').Value = '  -3.38%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('D32').Value = '''28.12'
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('D34').Value = '0.0₃0943'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '''48.92'
$ws.Range('E36').Value = '  +3.18%  '
$ws.Range('D37').Value = '''5.61'
$ws.Range('E37').Value = '  -4.55%  '
$ws.Range('D38').Value = '''0.948'
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').Value = '''49.19'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''1.99'
$ws.Range('E41').Value = '  -2.72%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.123'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').Value = '''8.32'
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('D44').Value = '''2.69'
$ws.Range('E44').Value = '  +3.27%  '
$ws.Range('D45').Value = '2.791.94'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').Value = '''373.00'
$ws.Range('E46').Value = '  -3.52%  '
$ws.Range('D47').Value = '''0.0345'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').Value = '''135.29'
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('D50').Value = '''25.04'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('E51').Value = '  +1.94%  '
